$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.334.10"
$ws.Range("E2").Value = "'  -4.16%  "
$ws.Range("D3").Value = "'2.617.04"
$ws.Range("E3").Value = "'  -3.96%  "
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'520.06"
$ws.Range("E5").Value = "'  -1.41%  "
$ws.Range("D6").Value = "'142.65"
$ws.Range("E6").Value = "'  -2.46%  "
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = "'  -1.68%  "
$ws.Range("D9").Value = "'6.57"
$ws.Range("E9").Value = "'  -5.81%  "
$ws.Range("E10").Value = "'  -2.69%  "
$ws.Range("D11").Value = "'0.337"
$ws.Range("E11").Value = "'  -0.77%  "
$ws.Range("E12").Value = "'  +1.03%  "
$ws.Range("D13").Value = "'3.078.24"
$ws.Range("E13").Value = "'  -3.19%  "
$ws.Range("D14").Value = "'58.312.47"
$ws.Range("E14").Value = "'  -4.10%  "
$ws.Range("E15").Value = "'  -1.94%  "
$ws.Range("E16").Value = "'  -1.59%  "
$ws.Range("D17").Value = "'2.625.63"
$ws.Range("E17").Value = "'  -7.42%  "
$ws.Range("D18").Value = "'336.60"
$ws.Range("E18").Value = "'  -2.42%  "
$ws.Range("E19").Value = "'  -3.02%  "
$ws.Range("E20").Value = "'  -1.55%  "
$ws.Range("E21").Value = "'  -2.62%  "
$ws.Range("E22").Value = "'  -0.09%  "
$ws.Range("D23").Value = "'64.40"
$ws.Range("E23").Value = "'  +1.03%  "
$ws.Range("E24").Value = "'  -1.21%  "
$ws.Range("E25").Value = "'  -1.74%  "
$ws.Range("E26").Value = "'  +0.35%  "
$ws.Range("D27").Value = "'7.12"
$ws.Range("E27").Value = "'  -2.32%  "
$ws.Range("D28").Value = "'0.0₃0789"
$ws.Range("E28").Value = "'  -3.95%  "
$ws.Range("D29").Value = "'6.55"
$ws.Range("E29").Value = "'  -3.89%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "'  +0.08%  "
$ws.Range("E31").Value = "'  -0.69%  "
$ws.Range("D32").Value = "'18.76"
$ws.Range("E32").Value = "'  -1.34%  "
$ws.Range("D33").Value = "'149.91"
$ws.Range("E33").Value = "'  +0.16%  "
$ws.Range("D34").Value = "'4.11"
$ws.Range("E34").Value = "'  -3.75%  "
$ws.Range("E35").Value = "'  -4.15%  "
$ws.Range("D36").Value = "'0.885"
$ws.Range("E36").Value = "'  -4.90%  "
$ws.Range("D37").Value = "'0.856"
$ws.Range("E37").Value = "'  -2.81%  "
$ws.Range("D38").Value = "'36.35"
$ws.Range("E38").Value = "'  -2.23%  "
$ws.Range("E39").Value = "'  -6.15%  "
$ws.Range("E40").Value = "'  -1.14%  "
$ws.Range("E41").Value = "'  -0.07%  "
$ws.Range("E42").Value = "'  -1.27%  "
$ws.Range("D43").Value = "'0.0968"
$ws.Range("E43").Value = "'  -1.93%  "
$ws.Range("D44").Value = "'269.03"
$ws.Range("E44").Value = "'  -4.59%  "
$ws.Range("E45").Value = "'  +1.04%  "
$ws.Range("D46").Value = "'19.13"
$ws.Range("E46").Value = "'  -5.29%  "
$ws.Range("D47").Value = "'0.0533"
$ws.Range("E47").Value = "'  -1.70%  "
$ws.Range("D48").Value = "'2.033.34"
$ws.Range("E48").Value = "'  -4.83%  "
$ws.Range("E49").Value = "'  -1.75%  "
$ws.Range("E50").Value = "'  -4.41%  "
$ws.Range("E51").Value = "'  -5.06%  "
